# Update "want-to-go" counts (column F) across the four sheets to reflect
# the latest scrape output (gh-pages generated data refresh).

$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F13").Value = 2700
$ws1.Range("F18").Value = 266
$ws1.Range("F20").Value = 5521
$ws1.Range("F25").Value = 412
$ws1.Range("F26").Value = 1182
$ws1.Range("F28").Value = 95
$ws1.Range("F29").Value = 311

# Sheet: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F13").Value = 626
$ws2.Range("F33").Value = 39

# Sheet: 本地生活
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 1791
$ws3.Range("F6").Value = 1095

# Sheet: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1791
$ws4.Range("F6").Value = 1095
$ws4.Range("F20").Value = 2700
$ws4.Range("F26").Value = 266
$ws4.Range("F28").Value = 5521
$ws4.Range("F31").Value = 626
$ws4.Range("F34").Value = 412
$ws4.Range("F41").Value = 1182
$ws4.Range("F47").Value = 39
$ws4.Range("F48").Value = 95
$ws4.Range("F49").Value = 311
